$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap columns B and C (header row + data row) per the new combined table method
$ws.Range("B1").Value = "total_requested_2020"
$ws.Range("C1").Value = "total_returned_2016"

$ws.Range("B2").Value = 61976781
$ws.Range("C2").Value = 21876963
